# Apply updated coin price/volume figures (and a CEJI/KickToken row swap)
# as scraped on Sun Jan  1 05:17:34 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E columns hold numbers-as-text (e.g. "243.94", "-0.59%"). A leading
# apostrophe forces Excel to keep the entry as text instead of
# auto-converting it to a real number/percentage; resetting the Style
# back to 'Normal' afterwards drops the quote-prefix formatting so the
# cell keeps its original (unstyled) look.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}
Set-TextValue 'D2' '243.94'
Set-TextValue 'E2' '-0.59%'
Set-TextValue 'D3' '26.29'
Set-TextValue 'E3' '3.38%'
Set-TextValue 'D4' '5.124'
Set-TextValue 'E4' '0.42%'
Set-TextValue 'D5' '0.05589'
Set-TextValue 'E5' '0.26%'
Set-TextValue 'D6' '6.465'
Set-TextValue 'E6' '-0.47%'
Set-TextValue 'D7' '0.8218'
Set-TextValue 'E7' '0.44%'
Set-TextValue 'D8' '0.8372'
Set-TextValue 'E8' '-1.07%'
Set-TextValue 'D9' '0.1332'
Set-TextValue 'E9' '-0.49%'
Set-TextValue 'D10' '0.06986'
Set-TextValue 'E10' '0.48%'
Set-TextValue 'D11' '0.02888'
Set-TextValue 'E11' '0.41%'
Set-TextValue 'D12' '0.09377'
Set-TextValue 'E12' '0.00%'
Set-TextValue 'E13' '0.44%'
Set-TextValue 'D14' '0.0005961'
Set-TextValue 'E14' '-0.15%'
Set-TextValue 'D15' '0.006223'
Set-TextValue 'E15' '2.05%'
Set-TextValue 'D16' '3.652'
Set-TextValue 'E16' '4.43%'
Set-TextValue 'D17' '3.032'
Set-TextValue 'E17' '0.44%'
Set-TextValue 'E18' '4.37%'
Set-TextValue 'E20' '-3.10%'
Set-TextValue 'E21' '-2.20%'
Set-TextValue 'D22' '3.745'
Set-TextValue 'E22' '-0.13%'
Set-TextValue 'D23' '0.04663'
Set-TextValue 'E23' '-1.30%'
Set-TextValue 'E24' '-0.05%'
Set-TextValue 'D26' '0.004500'
Set-TextValue 'E26' '-2.94%'
Set-TextValue 'D27' '0.00009601'
Set-TextValue 'E27' '-1.02%'
Set-TextValue 'D28' '0.0001393'
Set-TextValue 'E28' '0.27%'
Set-TextValue 'D40' '0.03640'
Set-TextValue 'E40' '-0.57%'
Set-TextValue 'D41' '0.1380'
Set-TextValue 'E41' '31.23%'
$ws.Range('B42').Value = 'KickToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue 'D42' '0.006151'
Set-TextValue 'E42' '-0.80%'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue 'D43' '0.002600'
Set-TextValue 'E43' '4.01%'
Set-TextValue 'D44' '0.008864'
Set-TextValue 'E44' '6.67%'
Set-TextValue 'D45' '0.00005355'
Set-TextValue 'E45' '1.13%'
Set-TextValue 'E46' '0.01%'
Set-TextValue 'E47' '-3.99%'
Set-TextValue 'E48' '10.05%'
Set-TextValue 'E49' '0.01%'
Set-TextValue 'E50' '0.01%'
